$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date-range banner in H1.
$ws.Range("H1").Value = "Bitacora de reportes del 26 de mayo de 2023 al 01 de junio de 2023"

# Rewrite row 9's log entry with the new shipment details
# (do this before deleting row 10 so the "9" row index still refers to it).
$ws.Range("A9").Value = 45078
$ws.Range("B9").Value = "Huerta Sierra Madre"
$ws.Range("C9").Value = "Pedro Mendoza Cobarrubias"
$ws.Range("D9").Value = "Aventajado"
$ws.Range("E9").Value = "Gasolinera Ejido Opopeo"
$ws.Range("F9").Value = "H-2345"
$ws.Range("G9").Value = "Ruben Juarez Hernandez"
$ws.Range("H9").Value = "Cuadrilla Michoacanos"
$ws.Range("I9").Value = "HINO12321"
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 10
$ws.Range("L9").Value = "Exportación"
$ws.Range("M9").Value = "Phoenix, Arizona"
$ws.Range("N9").Value = "Kevin Morales  Tellez"
$ws.Range("O9").Value = "El pedido debe enviar un regalo"

# Remove the old row 10 entry entirely; this shifts every row below it
# up by one (13->12, 14->13, 15->14, 18->17), matching the target layout.
$ws.Rows(10).Delete()
